$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.880.99"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "'2.643.30"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'581.24"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'155.85"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.622"
$ws.Range("E8").Value = "  -4.21%  "
$ws.Range("D9").Value = "'2.640.78"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -4.13%  "
$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'0.383"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "'28.46"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "'3.119.82"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "'63.790.02"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "'2.643.19"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'12.17"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "'7.72"
$ws.Range("D21").Value = "'4.54"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("D22").Value = "'345.90"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'68.06"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  +8.30%  "
$ws.Range("D26").Value = "'0.0000110"
$ws.Range("E26").Value = "  -3.84%  "
$ws.Range("D27").Value = "'609.12"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").Value = "'9.25"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("D29").Value = "'1.62"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'8.15"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'1.74"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'6.58"
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D39").Value = "'19.73"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").Value = "'1.90"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").Value = "'151.67"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").Value = "'41.95"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "'162.24"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "'3.91"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").Value = "'0.0589"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").Value = "'0.634"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").Value = "'0.0248"
$ws.Range("E51").Value = "  -3.20%  "
